$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParaXml($matchText, $innerXml) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $matchText) {
            $r = $d.Range($p.Range.Start, $p.Range.End - 1)
            $r.InsertXML("<w:p $wns>" + $innerXml + "</w:p>")
            return
        }
    }
}

# 1) "Satt opp GitHub" paragraph: merge runs, drop spell-check markers, add ".",
#    keep the trailing tab run (the _GoBack bookmark is relocated here afterwards).
Set-ParaXml "*Satt opp GitHub*" (
    "<w:r><w:t>Satt opp GitHub</w:t></w:r>" +
    "<w:r><w:t>.</w:t></w:r>" +
    "<w:r><w:tab/></w:r>"
)

# 2) "Laget en kjapp mockup" paragraph: drop spell-check markers around "mockup",
#    merge it into the first run, keep the remainder as its own run.
Set-ParaXml "*Laget en kjapp*" (
    "<w:r><w:t>Laget en kjapp mockup</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> av nettsiden som kan brukes som et utgangspunkt. Ligger vedlagt</w:t></w:r>"
)

# 3) "Hele gruppen jobber sammen via Git ..." paragraph: drop spell-check markers
#    around "Git" and merge all three runs into one.
Set-ParaXml "*Hele gruppen*" (
    "<w:r><w:t>Hele gruppen jobber sammen via Git for å ferdigstille dokumentasjonen.</w:t></w:r>"
)

# 4) "Begynner med å kode forside med header og footer." paragraph: drop spell-check
#    markers around "footer" and merge the three runs into one.
Set-ParaXml "*Begynner med*" (
    "<w:r><w:t>Begynner med å kode forside med header og footer.</w:t></w:r>"
)

# 5) "Jobber med hexagonmenyen ..." paragraph: drop spell-check markers around
#    "hexagonmenyen" and "responsiv" and merge all runs into one.
Set-ParaXml "*Jobber med hexagonmenyen*" (
    "<w:r><w:t>Jobber med hexagonmenyen og generell finpuss på diverse elementer. Prøver å lage den så responsiv jeg kan ut ifra kunnskaper fra forrige semester.</w:t></w:r>"
)

# 6) Move the _GoBack bookmark from the end of the document to right after the new
#    "." run in the "Satt opp GitHub." paragraph (this also removes it from its old
#    location, since a document only carries a single _GoBack bookmark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Satt opp GitHub*") {
        $pos = $p.Range.Start + "Satt opp GitHub.".Length
        $bmRange = $d.Range($pos, $pos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
